$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '63.530.58'
$ws.Range("E2").Value = '  +0.10%  '
$ws.Range("D3").Value = '2.581.64'
$ws.Range("E3").Value = '  -1.33%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '587.85'
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = '  -0.28%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '144.67'
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = '  -3.37%  '
$ws.Range("E7").Value = '  +0.17%  '
$ws.Range("E8").Value = '  -1.77%  '
$ws.Range("E9").Value = '  -3.70%  '
$ws.Range("E11").Value = '  -0.19%  '
$ws.Range("E12").Value = '  -1.99%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '27.31'
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = '  -2.23%  '
$ws.Range("D14").Value = '3.049.31'
$ws.Range("E14").Value = '  -1.05%  '
$ws.Range("D15").Value = '63.395.47'
$ws.Range("E15").Value = '  -0.06%  '
$ws.Range("E16").Value = '  -3.01%  '
$ws.Range("D17").Value = '2.583.26'
$ws.Range("E17").Value = '  -2.07%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '11.12'
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = '  -3.57%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '343.23'
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = '  -1.05%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '4.31'
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = '  -3.57%  '
$ws.Range("E21").Value = '  -3.98%  '
$ws.Range("E22").Value = '  +0.16%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '68.37'
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = '  +2.16%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '1.58'
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = '  +5.59%  '
$ws.Range("E25").Value = '  -1.67%  '
$ws.Range("E26").Value = '  -4.01%  '
$ws.Range("E27").Value = '  +0.23%  '
$ws.Range("E28").Value = '  -3.85%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '8.25'
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = '  -3.68%  '
$ws.Range("E30").Value = '  -2.54%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '470.84'
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = '  +0.20%  '
$ws.Range("E32").Value = '  -4.68%  '
$ws.Range("E33").Value = '  +1.45%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '176.62'
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = '  -0.37%  '
$ws.Range("E35").Value = '  +0.04%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.403'
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = '  -1.70%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '18.91'
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = '  -2.61%  '
$ws.Range("E38").Value = '  -3.43%  '
$ws.Range("E39").Value = '  -0.03%  '
$ws.Range("E40").Value = '  -2.16%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '162.09'
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = '  +5.30%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '40.04'
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = '  +0.90%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '3.72'
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = '  -3.89%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '21.78'
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = '  +2.57%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.630'
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = '  +1.78%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.0537'
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = '  -3.94%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.0961'
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = '  -1.98%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.0237'
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = '  -2.21%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '18.23'
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = '  -3.06%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '11.36'
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = '  -0.33%  '
$ws.Range("E51").Value = '  -4.95%  '
